$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date in column A (rows 2-28) from 1/12/2019 to 1/21/2019
$ws.Range("A2:A28").Value = "1/21/2019"

# Update the saved selection / active cell to C15
$ws.Range("C15").Select() | Out-Null
